$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2041522491349481
$ws.Range("C2").Value = 0.5674740484429066
$ws.Range("J2").Value = 0.02076124567474048
$ws.Range("P2").Value = 0.1384083044982699
$ws.Range("S2").Value = 0.06920415224913495
$ws.Range("C3").Value = 0.04705882352941176
$ws.Range("J3").Value = 0.01176470588235294
$ws.Range("O3").Value = 0.005882352941176471
$ws.Range("P3").Value = 0.7470588235294118
$ws.Range("S3").Value = 0.1882352941176471
$ws.Range("P4").Value = 0.7
$ws.Range("S4").Value = 0.3
$ws.Range("B6").Value = 0.06103286384976526
$ws.Range("E6").Value = 0.004694835680751174
$ws.Range("F6").Value = 0.07511737089201878
$ws.Range("J6").Value = 0.2206572769953052
$ws.Range("O6").Value = 0.03286384976525822
$ws.Range("Q6").Value = 0.1126760563380282
$ws.Range("R6").Value = 0.08450704225352113
$ws.Range("S6").Value = 0.4084507042253521
$ws.Range("B7").Value = 0.1016949152542373
$ws.Range("D7").Value = 0.02259887005649718
$ws.Range("E7").Value = 0.005649717514124294
$ws.Range("F7").Value = 0.05084745762711865
$ws.Range("J7").Value = 0.1016949152542373
$ws.Range("O7").Value = 0.005649717514124294
$ws.Range("Q7").Value = 0.1355932203389831
$ws.Range("R7").Value = 0.1016949152542373
$ws.Range("S7").Value = 0.4745762711864407
$ws.Range("B8").Value = 0.09073724007561437
$ws.Range("D8").Value = 0.01512287334593573
$ws.Range("F8").Value = 0.04158790170132325
$ws.Range("J8").Value = 0.0888468809073724
$ws.Range("O8").Value = 0.005671077504725898
$ws.Range("Q8").Value = 0.168241965973535
$ws.Range("R8").Value = 0.1190926275992439
$ws.Range("S8").Value = 0.4706994328922495
$ws.Range("B9").Value = 0.09375
$ws.Range("D9").Value = 0.01041666666666667
$ws.Range("F9").Value = 0.0625
$ws.Range("J9").Value = 0.125
$ws.Range("O9").Value = 0.005208333333333333
$ws.Range("Q9").Value = 0.140625
$ws.Range("R9").Value = 0.08854166666666667
$ws.Range("S9").Value = 0.4739583333333333
$ws.Range("B10").Value = 0.09720176730486009
$ws.Range("D10").Value = 0.01914580265095729
$ws.Range("E10").Value = 0.001472754050073638
$ws.Range("F10").Value = 0.07290132547864507
$ws.Range("J10").Value = 0.1060382916053019
$ws.Range("O10").Value = 0.009572901325478646
$ws.Range("Q10").Value = 0.1921944035346097
$ws.Range("R10").Value = 0.09351988217967599
$ws.Range("S10").Value = 0.4079528718703976
$ws.Range("G11").Value = 0.14453125
$ws.Range("J11").Value = 0.07421875
$ws.Range("K11").Value = 0.17578125
$ws.Range("L11").Value = 0.59375
$ws.Range("S11").Value = 0.01171875
$ws.Range("G12").Value = 0.7484276729559748
$ws.Range("J12").Value = 0.1823899371069182
$ws.Range("L12").Value = 0.03144654088050314
$ws.Range("S12").Value = 0.03773584905660377
$ws.Range("F13").Value = 0.02631578947368421
$ws.Range("G13").Value = 0.631578947368421
$ws.Range("J13").Value = 0.2105263157894737
$ws.Range("S13").Value = 0.131578947368421
$ws.Range("F15").Value = 0.01951219512195122
$ws.Range("H15").Value = 0.2
$ws.Range("I15").Value = 0.05853658536585366
$ws.Range("J15").Value = 0.4
$ws.Range("K15").Value = 0.05853658536585366
$ws.Range("M15").Value = 0.004878048780487805
$ws.Range("N15").Value = 0.004878048780487805
$ws.Range("O15").Value = 0.05365853658536585
$ws.Range("S15").Value = 0.2
$ws.Range("F16").Value = 0.01587301587301587
$ws.Range("H16").Value = 0.1957671957671958
$ws.Range("I16").Value = 0.04232804232804233
$ws.Range("J16").Value = 0.4391534391534391
$ws.Range("K16").Value = 0.1375661375661376
$ws.Range("M16").Value = 0.01587301587301587
$ws.Range("N16").Value = 0.005291005291005291
$ws.Range("O16").Value = 0.04232804232804233
$ws.Range("S16").Value = 0.1058201058201058
$ws.Range("F17").Value = 0.006976744186046512
$ws.Range("H17").Value = 0.2069767441860465
$ws.Range("I17").Value = 0.1046511627906977
$ws.Range("J17").Value = 0.4441860465116279
$ws.Range("K17").Value = 0.07674418604651163
$ws.Range("M17").Value = 0.004651162790697674
$ws.Range("N17").Value = 0.002325581395348837
$ws.Range("O17").Value = 0.04651162790697674
$ws.Range("S17").Value = 0.1069767441860465
$ws.Range("F18").Value = 0.008264462809917356
$ws.Range("H18").Value = 0.2231404958677686
$ws.Range("I18").Value = 0.09917355371900827
$ws.Range("J18").Value = 0.4173553719008264
$ws.Range("K18").Value = 0.07024793388429752
$ws.Range("M18").Value = 0.01652892561983471
$ws.Range("O18").Value = 0.07851239669421488
$ws.Range("S18").Value = 0.08677685950413223
$ws.Range("F19").Value = 0.01570306923625981
$ws.Range("H19").Value = 0.2177016416845111
$ws.Range("I19").Value = 0.07566024268379729
$ws.Range("J19").Value = 0.4154175588865097
$ws.Range("K19").Value = 0.08493932905067808
$ws.Range("M19").Value = 0.02069950035688794
$ws.Range("N19").Value = 0.0007137758743754461
$ws.Range("O19").Value = 0.06281227694503926
$ws.Range("S19").Value = 0.1063526052819415
